# Zenodo CHE_trd_coal.xlsx -> CHE_trd_gas.xlsx conversion
# Rename the "coal" trade commodity references to "gas" throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data rows (6-106) reference the parameter name "trd_coal" in column B;
# rename it to "trd_gas" for this fuel-specific workbook.
$ws.Range("B6:B106").Replace("trd_coal", "trd_gas") | Out-Null

# Row 16 holds the FxE (output_efficiency / constant_fxe) row whose "Flow" column (F)
# names the commodity itself; switch it from coal to gas.
$ws.Range("F16").Value = "gas"

# Move the active cell/selection to F17 (cosmetic, matches last user interaction).
$ws.Range("F17").Select() | Out-Null
